$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new "amazeloan" entry is logged on 2024-09-01. The existing most-recent
# entry (row 24) is pushed down by inserting a fresh row right below it, and
# that new row picks up the entry's original timestamp while row 24 itself
# is updated to the latest save time.
$ws.Rows(25).Insert()

# Row 24: refresh September_Date (S24) to the latest timestamp
$ws.Range("S24").Value = "2024-09-01 09:29:24"

# Row 25 (newly inserted): September details/date (R25/S25)
$ws.Range("R25").Value = "amazeloan"
$ws.Range("S25").Value = "2024-09-01 09:27:06"
